$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D price values that look like plain numbers (single decimal
# point, no thousands separators) get auto-converted from Text to Number by
# Excel COM when assigned via .Value, which silently drops formatting like
# trailing zeros (e.g. "105.60" -> 105.6). The source cells are plain text
# (t="inlineStr"), so prefix such values with a literal leading apostrophe
# (Excel's standard "treat as text" quote-prefix convention) to force them to
# stay Text and preserve the exact original string, then reset the cell style
# to Normal so the quote-prefix marker does not leave a stray style behind.

# Row 2
$ws.Range("D2").Value = "44.009.40"
$ws.Range("E2").Value = "  +3.01%  "

# Row 3
$ws.Range("D3").Value = "2.257.56"
$ws.Range("E3").Value = "  +2.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'258.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.71%  "

# Row 6
$ws.Range("D6").Value = "'79.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.92%  "

# Row 7
$ws.Range("E7").Value = "  +2.17%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.91%  "

# Row 10
$ws.Range("D10").Value = "'43.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.52%  "

# Row 11
$ws.Range("D11").Value = "'0.0929"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.15%  "

# Row 12
$ws.Range("D12").Value = "'7.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.28%  "

# Row 13
$ws.Range("E13").Value = "  +2.18%  "

# Row 14
$ws.Range("D14").Value = "2.590.44"
$ws.Range("E14").Value = "  +1.88%  "

# Row 15
$ws.Range("D15").Value = "'14.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.33%  "

# Row 16
$ws.Range("D16").Value = "2.238.97"
$ws.Range("E16").Value = "  +1.76%  "

# Row 17
$ws.Range("E17").Value = "  +2.13%  "

# Row 18
$ws.Range("D18").Value = "43.962.27"
$ws.Range("E18").Value = "  +3.15%  "

# Row 19
$ws.Range("E19").Value = "  +1.86%  "

# Row 20
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'71.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.75%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.58%  "

# Row 22
$ws.Range("D22").Value = "'2.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.98%  "

# Row 23
$ws.Range("D23").Value = "'235.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "

# Row 24
$ws.Range("D24").Value = "'9.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.15%  "

# Row 25
$ws.Range("E25").Value = "  -0.10%  "

# Row 26
$ws.Range("D26").Value = "'42.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.65%  "

# Row 27
$ws.Range("D27").Value = "'10.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.95%  "

# Row 28
$ws.Range("D28").Value = "'3.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "

# Row 29
$ws.Range("E29").Value = "  +1.48%  "

# Row 30
$ws.Range("E30").Value = "  -0.83%  "

# Row 31
$ws.Range("D31").Value = "'173.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.26%  "

# Row 32
$ws.Range("D32").Value = "'20.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.91%  "

# Row 33
$ws.Range("D33").Value = "'0.0880"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.11%  "

# Row 34
$ws.Range("D34").Value = "'5.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.46%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.87%  "

# Row 37
$ws.Range("D37").Value = "'4.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.84%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0362"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.36%  "

# Row 39
$ws.Range("D39").Value = "'13.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.55%  "

# Row 40
$ws.Range("D40").Value = "'2.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +19.45%  "

# Row 41
$ws.Range("E41").Value = "  +2.80%  "

# Row 42
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'62.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.00%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.205"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.08%  "

# Row 44
$ws.Range("E44").Value = "  +2.61%  "

# Row 45
$ws.Range("D45").Value = "'105.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.78%  "

# Row 46
$ws.Range("D46").Value = "'8.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.80%  "

# Row 47
$ws.Range("D47").Value = "'0.477"
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.0990"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "

# Row 49
$ws.Range("E49").Value = "  +2.48%  "

# Row 50
$ws.Range("E50").Value = "  +2.16%  "

# Row 51
$ws.Range("E51").Value = "  +26.28%  "
